# Commit: "updating timings based upon method inlining"
#
# The "Objeck (JIT)" benchmark times (column C, Sheet1) were re-measured
# after the interpreter started inlining methods, so they drop from
# ~0.87s to ~0.48s for each of the four recorded runs. C6 (=AVERAGE(C2:C5))
# and A8 (=ABS(C6/A6-1)) are formulas, so they recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$ws.Range("C2").Value = 0.47585300000000003
$ws.Range("C3").Value = 0.48208499999999999
$ws.Range("C4").Value = 0.47908600000000001
$ws.Range("C5").Value = 0.47655199999999998

# Move the active selection (as recorded in the saved view state) to I5.
$ws.Range("I5").Select()
